# Added Test Data for UK Market
#
# Adds a new "UK" worksheet (cloned from the last existing country sheet,
# "Poland", so it inherits the same layout/styles/merged cells/page setup),
# then overwrites the market-specific cells with the UK data, matching the
# structure used by every other country tab in this workbook.

$wb = $excel.ActiveWorkbook

# Clone the last sheet (Poland) to pick up formatting/merges/styles, and
# place the copy right after it - this becomes the new last tab, exactly
# like every previous "add a country" edit to this workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet) | Out-Null

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# Market code first, then market name - this mirrors shared-string
# allocation order from the source edit (B4's code string becomes the
# lower index, B2's market name the next one).
$newSheet.Range("B4").Value = "NGC-2741/T3365"
$newSheet.Range("B2").Value = "UK Market"

# Accessory list entries for the UK sheet (same two codes as Poland's
# sheet, but in the opposite row order).
$newSheet.Range("A9").Value = "MX-DPBX"
$newSheet.Range("A10").Value = "MX-BBX"

# New sheet becomes the active/selected tab, with B4 selected.
$newSheet.Range("B4").Select() | Out-Null

Write-Host "Added 'UK' worksheet with UK Market test data"
